# Refresh the "Metadata" sheet of the Bortezomib ValueSet workbook:
#  - bump Version 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date refreshed to the new publication timestamp
#  - Contact row split into the publisher contact + a named contact (Bob Milius)
#  - a new "Jurisdiction" row inserted right after the Contact rows, pushing
#    Description / Purpose / Copyright / Immutable down by one row and
#    adding a new trailing Immutable / BooleanType[null] row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value updates on the existing rows -----------------------------
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- give row 16 the same formatting as the rest of the data rows (copy the
#     formats only - PasteSpecial xlPasteFormats keeps the existing cellXfs
#     table intact, unlike Rows.Insert which mints a fresh, unused style) --
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# --- rewrite rows 12-16 with their final (post-shift) content -------------
# row 12: new "Jurisdiction" property (no value)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
# row 13: was row 12 (Description)
$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "RxNorm codes for Bortezomib"
# row 14: was row 13 (Purpose, no value)
$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").ClearContents()
# row 15: was row 14 (Copyright, no value)
$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").ClearContents()
# row 16: new trailing row, was row 15 (Immutable / BooleanType[null])
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
